$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A/D labels for rows 4 and 5
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("D4").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"

# Row 2 values (E2:T2)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3642143333333334
$ws.Range("H2").Value = 1.092643
$ws.Range("I2").Value = 0.4800482050304226
$ws.Range("J2").Value = 0.4800482050304224
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.467406
$ws.Range("N2").Value = 4.402218
$ws.Range("O2").Value = 0.864087546066766
$ws.Range("P2").Value = 0.9050919696083439
$ws.Range("Q2").Value = 0.5344502980193333
$ws.Range("R2").Value = 4.810052682174
$ws.Range("S2").Value = 0.4148036754784936
$ws.Range("T2").Value = 0.4344877753979351

# Row 3 values (E3:T3)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3642143333333334
$ws.Range("H3").Value = 1.092643
$ws.Range("I3").Value = 0.4800482050304226
$ws.Range("J3").Value = 0.4800482050304224
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2308085
$ws.Range("N3").Value = 0.461617
$ws.Range("O3").Value = 0.135912453933234
$ws.Range("P3").Value = 0.09490803039165596
$ws.Range("Q3").Value = 0.08406376395516667
$ws.Range("R3").Value = 0.504382583731
$ws.Range("S3").Value = 0.06524452955192898
$ws.Range("T3").Value = 0.04556042963248723

# Row 4 values (E4:T4)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3944893333333333
$ws.Range("H4").Value = 1.183468
$ws.Range("I4").Value = 0.5199517949695774
$ws.Range("J4").Value = 0.5199517949695774
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.467406
$ws.Range("N4").Value = 4.402218
$ws.Range("O4").Value = 0.864087546066766
$ws.Range("P4").Value = 0.9050919696083439
$ws.Range("Q4").Value = 0.5788760146693332
$ws.Range("R4").Value = 5.209884132023999
$ws.Range("S4").Value = 0.4492838705882724
$ws.Range("T4").Value = 0.4706041942104087

# Row 5 values (E5:T5)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3944893333333333
$ws.Range("H5").Value = 1.183468
$ws.Range("I5").Value = 0.5199517949695774
$ws.Range("J5").Value = 0.5199517949695774
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2308085
$ws.Range("N5").Value = 0.461617
$ws.Range("O5").Value = 0.135912453933234
$ws.Range("P5").Value = 0.09490803039165596
$ws.Range("Q5").Value = 0.09105149129266665
$ws.Range("R5").Value = 0.546308947756
$ws.Range("S5").Value = 0.07066792438130502
$ws.Range("T5").Value = 0.04934760075916873

# Delete rows 6 and 7 (old extra rows)
$ws.Range("A6:T7").Delete()
